# Finished input list generation and sanity checks:
# regenerate trial_total (F) sequentially from 297, refresh the reshuffled
# stimulus/condition columns (I-S) for each trial row, and fix the
# catch-trial image reference (L12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 297
$ws.Range("L2").Value = "stimuli/img_pna7l.png"
$ws.Range("M2").Value = 85.53333333333333
$ws.Range("N2").Value = 67.97777777777777
$ws.Range("O2").Value = 76.75555555555556
$ws.Range("P2").Value = 45
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 9
$ws.Range("F3").Value = 298
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = "new"
$ws.Range("K3").Value = "f"
$ws.Range("L3").Value = "stimuli/img_pdzf1.png"
$ws.Range("M3").Value = 86.23913043478261
$ws.Range("N3").Value = 67.17391304347827
$ws.Range("O3").Value = 76.70652173913044
$ws.Range("Q3").Value = 9
$ws.Range("R3").Value = 9
$ws.Range("S3").Value = 9
$ws.Range("F4").Value = 299
$ws.Range("I4").Value = $null
$ws.Range("J4").Value = "new"
$ws.Range("K4").Value = "f"
$ws.Range("L4").Value = "stimuli/img_3jnt7.png"
$ws.Range("M4").Value = 49.52272727272727
$ws.Range("N4").Value = 35.25
$ws.Range("O4").Value = 42.38636363636364
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 3
$ws.Range("F5").Value = 300
$ws.Range("I5").Value = "target"
$ws.Range("J5").Value = "old"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_9oofc.png"
$ws.Range("M5").Value = 82.47619047619048
$ws.Range("N5").Value = 65.5
$ws.Range("O5").Value = 73.98809523809524
$ws.Range("P5").Value = 42
$ws.Range("Q5").Value = 8
$ws.Range("R5").Value = 8
$ws.Range("S5").Value = 8
$ws.Range("F6").Value = 301
$ws.Range("L6").Value = "stimuli/img_2qhro.png"
$ws.Range("M6").Value = 81.73809523809524
$ws.Range("N6").Value = 62.73809523809524
$ws.Range("O6").Value = 72.23809523809524
$ws.Range("P6").Value = 42
$ws.Range("Q6").Value = 8
$ws.Range("R6").Value = 8
$ws.Range("S6").Value = 8
$ws.Range("F7").Value = 302
$ws.Range("L7").Value = "stimuli/img_amsgw.png"
$ws.Range("M7").Value = 86.08510638297872
$ws.Range("N7").Value = 65.95744680851064
$ws.Range("O7").Value = 76.02127659574468
$ws.Range("P7").Value = 47
$ws.Range("Q7").Value = 9
$ws.Range("R7").Value = 9
$ws.Range("S7").Value = 9
$ws.Range("F8").Value = 303
$ws.Range("L8").Value = "stimuli/img_iudc4.png"
$ws.Range("M8").Value = 73.625
$ws.Range("N8").Value = 52.275
$ws.Range("O8").Value = 62.95
$ws.Range("P8").Value = 40
$ws.Range("Q8").Value = 6
$ws.Range("R8").Value = 6
$ws.Range("S8").Value = 6
$ws.Range("F9").Value = 304
$ws.Range("L9").Value = "stimuli/img_rg4in.png"
$ws.Range("M9").Value = 49.3695652173913
$ws.Range("N9").Value = 30.21739130434782
$ws.Range("O9").Value = 39.79347826086956
$ws.Range("P9").Value = 46
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 3
$ws.Range("S9").Value = 3
$ws.Range("F10").Value = 305
$ws.Range("L10").Value = "stimuli/img_xr3up.png"
$ws.Range("M10").Value = 76.24444444444444
$ws.Range("N10").Value = 55.88888888888889
$ws.Range("O10").Value = 66.06666666666666
$ws.Range("P10").Value = 45
$ws.Range("Q10").Value = 7
$ws.Range("R10").Value = 7
$ws.Range("S10").Value = 7
$ws.Range("F11").Value = 306
$ws.Range("L11").Value = "stimuli/img_zxvl3.png"
$ws.Range("M11").Value = 68.78260869565217
$ws.Range("N11").Value = 47.56521739130435
$ws.Range("O11").Value = 58.17391304347827
$ws.Range("P11").Value = 46
$ws.Range("Q11").Value = 5
$ws.Range("R11").Value = 5
$ws.Range("S11").Value = 5
$ws.Range("F12").Value = 307
$ws.Range("L12").Value = "stimuli/catch_10.jpg"
$ws.Range("F13").Value = 308
$ws.Range("I13").Value = "target"
$ws.Range("J13").Value = "old"
$ws.Range("K13").Value = "j"
$ws.Range("L13").Value = "stimuli/img_of8d6.png"
$ws.Range("M13").Value = 26.04878048780488
$ws.Range("N13").Value = 19.14634146341463
$ws.Range("O13").Value = 22.59756097560975
$ws.Range("P13").Value = 41
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 1
$ws.Range("S13").Value = 1
$ws.Range("F14").Value = 309
$ws.Range("I14").Value = "target"
$ws.Range("J14").Value = "old"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_vh7v8.png"
$ws.Range("M14").Value = 78.70454545454545
$ws.Range("N14").Value = 59.63636363636363
$ws.Range("O14").Value = 69.17045454545455
$ws.Range("P14").Value = 44
$ws.Range("Q14").Value = 7
$ws.Range("R14").Value = 7
$ws.Range("S14").Value = 7
$ws.Range("F15").Value = 310
$ws.Range("L15").Value = "stimuli/img_ra2nm.png"
$ws.Range("M15").Value = 70.75
$ws.Range("N15").Value = 50.375
$ws.Range("O15").Value = 60.5625
$ws.Range("P15").Value = 40
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6
$ws.Range("F16").Value = 311
$ws.Range("L16").Value = "stimuli/img_j4ttn.png"
$ws.Range("M16").Value = 12.61904761904762
$ws.Range("N16").Value = 11.42857142857143
$ws.Range("O16").Value = 12.02380952380952
$ws.Range("P16").Value = 42
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1
$ws.Range("F17").Value = 312
$ws.Range("L17").Value = "stimuli/img_syam3.png"
$ws.Range("M17").Value = 41.32432432432432
$ws.Range("N17").Value = 26.2972972972973
$ws.Range("O17").Value = 33.81081081081081
$ws.Range("P17").Value = 37
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 2
$ws.Range("F18").Value = 313
$ws.Range("L18").Value = "stimuli/img_qdln8.png"
$ws.Range("M18").Value = 85.51162790697674
$ws.Range("N18").Value = 67.86046511627907
$ws.Range("O18").Value = 76.68604651162791
$ws.Range("P18").Value = 43
$ws.Range("F19").Value = 314
$ws.Range("I19").Value = $null
$ws.Range("J19").Value = "new"
$ws.Range("K19").Value = "f"
$ws.Range("L19").Value = "stimuli/img_lgxzn.png"
$ws.Range("M19").Value = 73.11363636363636
$ws.Range("N19").Value = 49.97727272727273
$ws.Range("O19").Value = 61.54545454545455
$ws.Range("P19").Value = 44
$ws.Range("Q19").Value = 6
$ws.Range("R19").Value = 6
$ws.Range("S19").Value = 6
$ws.Range("F20").Value = 315
$ws.Range("L20").Value = "stimuli/img_rru0v.png"
$ws.Range("M20").Value = 56.45238095238095
$ws.Range("N20").Value = 39.42857142857143
$ws.Range("O20").Value = 47.94047619047619
$ws.Range("P20").Value = 42
$ws.Range("Q20").Value = 4
$ws.Range("R20").Value = 4
$ws.Range("S20").Value = 4
$ws.Range("F21").Value = 316
$ws.Range("L21").Value = "stimuli/img_24rt2.png"
$ws.Range("M21").Value = 55.26829268292683
$ws.Range("N21").Value = 34.19512195121951
$ws.Range("O21").Value = 44.73170731707317
$ws.Range("P21").Value = 41
$ws.Range("Q21").Value = 3
$ws.Range("R21").Value = 3
$ws.Range("S21").Value = 3
$ws.Range("F22").Value = 317
$ws.Range("L22").Value = "stimuli/img_ac0ey.png"
$ws.Range("M22").Value = 86.62222222222222
$ws.Range("N22").Value = 70.02222222222223
$ws.Range("O22").Value = 78.32222222222222
$ws.Range("P22").Value = 45
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 9
$ws.Range("F23").Value = 318
$ws.Range("I23").Value = "target"
$ws.Range("J23").Value = "old"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_vgh2g.png"
$ws.Range("M23").Value = 93.81395348837209
$ws.Range("N23").Value = 78.27906976744185
$ws.Range("O23").Value = 86.04651162790697
$ws.Range("P23").Value = 43
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = 10
$ws.Range("S23").Value = 10
$ws.Range("F24").Value = 319
$ws.Range("L24").Value = "stimuli/img_kq9s9.png"
$ws.Range("M24").Value = 62.30232558139535
$ws.Range("N24").Value = 39.97674418604651
$ws.Range("O24").Value = 51.13953488372093
$ws.Range("P24").Value = 43
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = 4
$ws.Range("S24").Value = 4
$ws.Range("F25").Value = 320
$ws.Range("I25").Value = $null
$ws.Range("J25").Value = "new"
$ws.Range("K25").Value = "f"
$ws.Range("L25").Value = "stimuli/img_7lz7m.png"
$ws.Range("M25").Value = 51.5531914893617
$ws.Range("N25").Value = 32.87234042553192
$ws.Range("O25").Value = 42.21276595744681
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = 3
$ws.Range("S25").Value = 3
$ws.Range("F26").Value = 321
$ws.Range("L26").Value = "stimuli/img_tn8ys.png"
$ws.Range("M26").Value = 86.70454545454545
$ws.Range("N26").Value = 72.4090909090909
$ws.Range("O26").Value = 79.55681818181819
$ws.Range("P26").Value = 44
$ws.Range("Q26").Value = 10
$ws.Range("R26").Value = 10
$ws.Range("S26").Value = 10
$ws.Range("F27").Value = 322
$ws.Range("L27").Value = "stimuli/img_swq34.png"
$ws.Range("M27").Value = 64.11363636363636
$ws.Range("N27").Value = 43.04545454545455
$ws.Range("O27").Value = 53.57954545454545
$ws.Range("P27").Value = 44
$ws.Range("Q27").Value = 5
$ws.Range("R27").Value = 5
$ws.Range("S27").Value = 5
$ws.Range("F28").Value = 323
$ws.Range("L28").Value = "stimuli/img_bf8nx.png"
$ws.Range("M28").Value = 86.63414634146342
$ws.Range("N28").Value = 66.63414634146342
$ws.Range("O28").Value = 76.63414634146342
$ws.Range("P28").Value = 41
$ws.Range("F29").Value = 324
$ws.Range("L29").Value = "stimuli/img_rych7.png"
$ws.Range("M29").Value = 30.4468085106383
$ws.Range("N29").Value = 23.4468085106383
$ws.Range("O29").Value = 26.9468085106383
$ws.Range("P29").Value = 47
$ws.Range("Q29").Value = 2
$ws.Range("R29").Value = 2
$ws.Range("S29").Value = 2
